$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10, pushing existing rows 10-14 down to 11-15
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly data
$ws.Cells.Item(10, 1).Value = 12
$ws.Cells.Item(10, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44477
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100112028
$ws.Cells.Item(10, 7).Value = "Sandia"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 800
$ws.Cells.Item(10, 12).Value = 800
$ws.Cells.Item(10, 13).Value = 800
$ws.Cells.Item(10, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 800
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# Copy the date-cell style (numFmt) from row 9's date cell into the new row's date cell
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4122)
$ws.Cells.Item(10, 4).Value = 44477
